$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data: product name in column A, unit in column B
$ws.Range("A2").Value = "Fresh Milk Green Field"
$ws.Range("B2").Value = "ml"
$ws.Range("A3").Value = "UHT Milk Indomilk"
$ws.Range("B3").Value = "ml"

# Match column A style to the header row's A1 style
$ws.Range("A2:A3").Font.Name = $ws.Range("A1").Font.Name
$ws.Range("A2:A3").Font.Size = $ws.Range("A1").Font.Size
$ws.Range("A2:A3").WrapText = $ws.Range("A1").WrapText
$ws.Range("A2:A3").VerticalAlignment = $ws.Range("A1").VerticalAlignment

# Column B (units) centered, using Calibri 11 font, horizontal+vertical center/bottom
$ws.Range("B2:B3").Font.Name = "Calibri"
$ws.Range("B2:B3").Font.Size = 11
$ws.Range("B2:B3").HorizontalAlignment = -4108
$ws.Range("B2:B3").VerticalAlignment = -4107

# Also update B1 to the same style as B2/B3 per diff (font id 2 removed, style now uses font 1 + centered)
$ws.Range("B1").Font.Name = "Calibri"
$ws.Range("B1").Font.Size = 11
$ws.Range("B1").HorizontalAlignment = -4108
$ws.Range("B1").VerticalAlignment = -4107
